$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'68.372.69"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  +1.33%  '
$ws.Range("D3").Value = "'3.743.78"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -0.45%  '
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.01%  '
$ws.Range("D5").Value = "'594.50"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.17%  '
$ws.Range("D6").Value = "'166.85"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -0.78%  '
$ws.Range("D7").Value = "'3.740.67"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.47%  '
$ws.Range("E9").Value = '  -0.79%  '
$ws.Range("E10").Value = '  -2.67%  '
$ws.Range("D11").Value = "'6.48"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.29%  '
$ws.Range("E12").Value = '  -1.16%  '
$ws.Range("E13").Value = '  -5.17%  '
$ws.Range("D14").Value = "'36.06"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  -0.73%  '
$ws.Range("D15").Value = "'4.366.71"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.44%  '
$ws.Range("D16").Value = "'3.744.34"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.52%  '
$ws.Range("D17").Value = "'68.354.53"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.39%  '
$ws.Range("D18").Value = "'17.86"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -3.63%  '
$ws.Range("D19").Value = "'7.01"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -2.09%  '
$ws.Range("E20").Value = '  -0.15%  '
$ws.Range("D21").Value = "'10.71"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +1.91%  '
$ws.Range("D22").Value = "'466.38"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.09%  '
$ws.Range("D23").Value = "'0.699"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -2.32%  '
$ws.Range("D24").Value = "'83.94"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +0.54%  '
$ws.Range("E25").Value = '  -1.33%  '
$ws.Range("D26").Value = "'2.19"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.08%  '
$ws.Range("D27").Value = "'12.01"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.76%  '
$ws.Range("D28").Value = "'10.13"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.79%  '
$ws.Range("E29").Value = '  -0.11%  '
$ws.Range("D30").Value = "'3.886.46"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -0.52%  '
$ws.Range("D31").Value = "'2.77"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -4.38%  '
$ws.Range("E32").Value = '  -4.26%  '
$ws.Range("D33").Value = "'29.85"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.72%  '
$ws.Range("E34").Value = '  -1.73%  '
$ws.Range("D35").Value = "'9.22"
$ws.Range("D35").Style = "Normal"
$ws.Range("D37").Value = "'3.695.42"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.62%  '
$ws.Range("E38").Value = '  -2.21%  '
$ws.Range("D39").Value = "'3.39"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -10.52%  '
$ws.Range("E40").Value = '  +0.16%  '
$ws.Range("D41").Value = "'0.995"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.42%  '
$ws.Range("D42").Value = "'5.78"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +0.12%  '
$ws.Range("E43").Value = '  +0.02%  '
$ws.Range("E45").Value = '  -1.86%  '
$ws.Range("D46").Value = "'8.58"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.83%  '
$ws.Range("E47").Value = '  -0.47%  '
$ws.Range("D48").Value = "'42.84"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +10.40%  '
$ws.Range("D49").Value = "'45.83"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").Value = "'146.12"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +5.07%  '
$ws.Range("D51").Value = "'390.36"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -1.21%  '
